# Apply price/volume refresh values scraped on Sat Sep 30 08:16:41 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as text, matching the
# original inline-string cell content (avoids numeric auto-conversion / rounding).
$ws.Range("D2").Value = "'26.954.54"
$ws.Range("E2").Value = "'  -0.65%  "
$ws.Range("D3").Value = "'1.674.04"
$ws.Range("E3").Value = "'  -0.34%  "
$ws.Range("E4").Value = "'  +0.15%  "
$ws.Range("D5").Value = "'214.69"
$ws.Range("E5").Value = "'  -0.95%  "
$ws.Range("E6").Value = "'  +1.55%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("E8").Value = "'  -0.98%  "
$ws.Range("E9").Value = "'  -0.24%  "
$ws.Range("D10").Value = "'20.36"
$ws.Range("E10").Value = "'  +0.51%  "
$ws.Range("D11").Value = "'0.0885"
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("D12").Value = "'1.910.28"
$ws.Range("E12").Value = "'  -0.37%  "
$ws.Range("D13").Value = "'1.645.43"
$ws.Range("E13").Value = "'  -2.09%  "
$ws.Range("E14").Value = "'  -0.60%  "
$ws.Range("E15").Value = "'  +0.24%  "
$ws.Range("D16").Value = "'65.65"
$ws.Range("E16").Value = "'  -0.86%  "
$ws.Range("D17").Value = "'26.963.78"
$ws.Range("D18").Value = "'236.12"
$ws.Range("E18").Value = "'  -1.46%  "
$ws.Range("D19").Value = "'8.10"
$ws.Range("E19").Value = "'  +4.19%  "
$ws.Range("E20").Value = "'  -1.06%  "
$ws.Range("E21").Value = "'  +0.14%  "
$ws.Range("E22").Value = "'  -1.38%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "'  -1.40%  "
$ws.Range("E24").Value = "'  -2.35%  "
$ws.Range("D25").Value = "'145.40"
$ws.Range("E25").Value = "'  -0.27%  "
$ws.Range("D26").Value = "'7.23"
$ws.Range("E26").Value = "'  +0.92%  "
$ws.Range("D27").Value = "'16.02"
$ws.Range("E27").Value = "'  -0.03%  "
$ws.Range("E28").Value = "'  -1.52%  "
$ws.Range("E29").Value = "'  +0.21%  "
$ws.Range("E30").Value = "'  -0.52%  "
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("D32").Value = "'3.31"
$ws.Range("E32").Value = "'  -0.21%  "
$ws.Range("D33").Value = "'1.479.29"
$ws.Range("E33").Value = "'  -0.76%  "
$ws.Range("E34").Value = "'  -0.24%  "
$ws.Range("D35").Value = "'1.67"
$ws.Range("E35").Value = "'  +3.01%  "
$ws.Range("E36").Value = "'  +0.11%  "
$ws.Range("D37").Value = "'0.587"
$ws.Range("E37").Value = "'  +1.24%  "
$ws.Range("D38").Value = "'0.895"
$ws.Range("E38").Value = "'  -1.29%  "
$ws.Range("D39").Value = "'0.0171"
$ws.Range("E39").Value = "'  +0.41%  "
$ws.Range("D40").Value = "'5.89"
$ws.Range("E40").Value = "'  -2.88%  "
$ws.Range("E41").Value = "'  +5.01%  "
$ws.Range("E42").Value = "'  +0.17%  "
$ws.Range("D43").Value = "'2.31"
$ws.Range("D44").Value = "'66.89"
$ws.Range("E44").Value = "'  -0.45%  "
$ws.Range("D45").Value = "'1.814.69"
$ws.Range("E45").Value = "'  -0.62%  "
$ws.Range("E46").Value = "'  -0.48%  "
$ws.Range("D47").Value = "'90.53"
$ws.Range("E47").Value = "'  +0.02%  "
$ws.Range("E48").Value = "'  -0.70%  "
$ws.Range("E49").Value = "'  +0.55%  "
$ws.Range("D50").Value = "'0.0507"
$ws.Range("E50").Value = "'  +0.01%  "
$ws.Range("D51").Value = "'7.65"
$ws.Range("E51").Value = "'  -0.63%  "
